$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "2022" data column (Z) mirroring the existing "2021" column (Y)
# Copy the formatting of each Y cell down to the corresponding Z cell first,
# so the new column inherits identical number formats / fonts / borders,
# then write in the new values.
# ---------------------------------------------------------------------------
$ws.Range("Y4").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("Z4").Value2 = 2022

$zValues = @{
    5  = 47.345690436648667
    6  = 55.294335329978139
    7  = 42.721146742902135
    8  = 56.732662465911261
    9  = 39.351829932862628
    10 = 43.952035422218046
    11 = 57.461907794486649
    12 = 32.073481974524846
    13 = 33.564455947162017
    14 = 55.803694659011171
    15 = 63.920911723512503
    16 = 52.521342498654128
}

foreach ($r in 5..16) {
    $ws.Range("Y$r").Copy()
    $ws.Range("Z$r").PasteSpecial(-4122)
    $ws.Range("Z$r").Value2 = $zValues[$r]
}

# ---------------------------------------------------------------------------
# Header row (row 4) re-alignment: the "Items" header cells now center their
# text both horizontally and vertically.
# ---------------------------------------------------------------------------
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108

$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Move the active selection to reflect the extra column of data (AA4 instead
# of AA15).
# ---------------------------------------------------------------------------
$ws.Range("AA4").Select()
